# Apply HotStock_Top20 update (2025-12-23 auto-generated refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "英维克"
$ws.Range("B2").Value = "平潭发展"
$ws.Range("C2").Value = "平潭发展"
$ws.Range("A3").Value = "海南发展"
$ws.Range("B3").Value = "东百集团"
$ws.Range("C3").Value = "东百集团"
$ws.Range("A4").Value = "平潭发展"
$ws.Range("B4").Value = "海南发展"
$ws.Range("C4").Value = "航天发展"
$ws.Range("A5").Value = "东百集团"
$ws.Range("B5").Value = "多氟多"
$ws.Range("C5").Value = "纳百川"
$ws.Range("A6").Value = "航天发展"
$ws.Range("B6").Value = "山子高科"
$ws.Range("C6").Value = "永辉超市"
$ws.Range("A7").Value = "神剑股份"
$ws.Range("B7").Value = "神剑股份"
$ws.Range("C7").Value = "神剑股份"
$ws.Range("A8").Value = "山子高科"
$ws.Range("B8").Value = "中国中免"
$ws.Range("C8").Value = "英维克"
$ws.Range("A9").Value = "永辉超市"
$ws.Range("B9").Value = "航天发展"
$ws.Range("C9").Value = "锡华科技"
$ws.Range("A10").Value = "N纳百川"
$ws.Range("B10").Value = "英维克"
$ws.Range("C10").Value = "再升科技"
$ws.Range("A11").Value = "中国中免"
$ws.Range("B11").Value = "海南瑞泽"
$ws.Range("C11").Value = "保变电气"
$ws.Range("A12").Value = "多氟多"
$ws.Range("B12").Value = "天际股份"
$ws.Range("C12").Value = "浙江世宝"
$ws.Range("A13").Value = "中国卫星"
$ws.Range("B13").Value = "海南海药"
$ws.Range("C13").Value = "博纳影业"
$ws.Range("A14").Value = "海南海药"
$ws.Range("B14").Value = "永辉超市"
$ws.Range("C14").Value = "海南发展"
$ws.Range("A15").Value = "天际股份"
$ws.Range("B15").Value = "N纳百川"
$ws.Range("C15").Value = "鹭燕医药"
$ws.Range("A16").Value = "N锡华"
$ws.Range("B16").Value = "中国卫星"
$ws.Range("C16").Value = "山子高科"
$ws.Range("A17").Value = "圣晖集成"
$ws.Range("B17").Value = "天赐材料"
$ws.Range("C17").Value = "航天机电"
$ws.Range("A18").Value = "中百集团"
$ws.Range("B18").Value = "中百集团"
$ws.Range("C18").Value = "中百集团"
$ws.Range("A19").Value = "亚翔集成"
$ws.Range("B19").Value = "航天电子"
$ws.Range("C19").Value = "安记食品"
$ws.Range("A20").Value = "浙江世宝"
$ws.Range("B20").Value = "神农种业"
$ws.Range("C20").Value = "西部材料"
$ws.Range("A21").Value = "海南瑞泽"
$ws.Range("B21").Value = "N锡华"
$ws.Range("C21").Value = "雪人集团"
